# The first paragraph in the formula document contains the first
# embedded equation (an Equation.KSEE3 OLE object run) followed by a
# trailing run of four plain spaces, and then the _GoBack bookmark.
# That trailing whitespace-only run is stray leftover text that should
# be removed, leaving just the equation run followed by the bookmark.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$paraRange = $p1.Range

# Exclude the paragraph mark (last character) from the range so only
# the visible run content -- the four trailing spaces -- is targeted.
$trailing = $d.Range($paraRange.Start, $paraRange.End - 1)

if ($trailing.Text -eq "    ") {
    $trailing.Text = ""
}
